$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# 1. Duplicate the "总计" sheet, placing the copy right after it.
$totalSheet.Copy($null, $totalSheet)

# 2. Rename: original becomes "2022-Q1" (new quarter detail data),
#    the fresh copy becomes the new "总计" (summary) sheet.
$totalSheet.Name = "2022-Q1"
$newTotal = $wb.Worksheets.Item("总计 (2)")
$newTotal.Name = "总计"

$q1 = $wb.Worksheets.Item("2022-Q1")

# --- Rebuild the "2022-Q1" sheet with fund-holding detail rows ---
# Extend header style (copied from the existing D1 header cell) across E1:H1.
$q1.Range("D1").Copy($q1.Range("E1:H1"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Extend the index-column style (copied from A2) down through A7:A11 for the new rows.
$q1.Range("A2").Copy($q1.Range("A7:A11"))

# Columns B (fund code) and D:G (scale/position/value) must stay TEXT, not be
# coerced to numbers (leading zeros in fund codes, decimal text in D:G) -- so
# pre-format those columns as Text before writing the numeric-looking strings.
$q1.Range("B2:B11").NumberFormat = "@"
$q1.Range("D2:G11").NumberFormat = "@"

$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "001583"
$q1.Cells.Item(2, 3).Value = "安信新常态沪港深精选股票"
$q1.Cells.Item(2, 4).Value = "8.85"
$q1.Cells.Item(2, 5).Value = "89.85"
$q1.Cells.Item(2, 6).Value = "3.59"
$q1.Cells.Item(2, 7).Value = "0.3177"
$q1.Cells.Item(2, 8).Value = 7

$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "011905"
$q1.Cells.Item(3, 3).Value = "安信价值启航混合型证券投资基金A"
$q1.Cells.Item(3, 4).Value = "4.48"
$q1.Cells.Item(3, 5).Value = "93.14"
$q1.Cells.Item(3, 6).Value = "3.21"
$q1.Cells.Item(3, 7).Value = "0.1438"
$q1.Cells.Item(3, 8).Value = 7

$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).Value = "011355"
$q1.Cells.Item(4, 3).Value = "华泰柏瑞港股通时代机遇混合型证券投资基金A"
$q1.Cells.Item(4, 4).Value = "1.13"
$q1.Cells.Item(4, 5).Value = "90.93"
$q1.Cells.Item(4, 6).Value = "8.07"
$q1.Cells.Item(4, 7).Value = "0.0912"
$q1.Cells.Item(4, 8).Value = 3

$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).Value = "008477"
$q1.Cells.Item(5, 3).Value = "安信价值驱动三年持有期混合"
$q1.Cells.Item(5, 4).Value = "2.64"
$q1.Cells.Item(5, 5).Value = "91.75"
$q1.Cells.Item(5, 6).Value = "3.40"
$q1.Cells.Item(5, 7).Value = "0.0898"
$q1.Cells.Item(5, 8).Value = 7

$q1.Cells.Item(6, 1).Value = 4
$q1.Cells.Item(6, 2).Value = "007109"
$q1.Cells.Item(6, 3).Value = "南方沪港深核心优势混合"
$q1.Cells.Item(6, 4).Value = "1.82"
$q1.Cells.Item(6, 5).Value = "87.54"
$q1.Cells.Item(6, 6).Value = "3.23"
$q1.Cells.Item(6, 7).Value = "0.0588"
$q1.Cells.Item(6, 8).Value = 6

$q1.Cells.Item(7, 1).Value = 5
$q1.Cells.Item(7, 2).Value = "161229"
$q1.Cells.Item(7, 3).Value = "国投瑞银中国价值发现股票QDII-LOF"
$q1.Cells.Item(7, 4).Value = "1.47"
$q1.Cells.Item(7, 5).Value = "92.83"
$q1.Cells.Item(7, 6).Value = "3.62"
$q1.Cells.Item(7, 7).Value = "0.0532"
$q1.Cells.Item(7, 8).Value = 10

$q1.Cells.Item(8, 1).Value = 6
$q1.Cells.Item(8, 2).Value = "003413"
$q1.Cells.Item(8, 3).Value = "华泰柏瑞新经济沪港深灵活配置混合"
$q1.Cells.Item(8, 4).Value = "0.54"
$q1.Cells.Item(8, 5).Value = "92.57"
$q1.Cells.Item(8, 6).Value = "8.82"
$q1.Cells.Item(8, 7).Value = "0.0476"
$q1.Cells.Item(8, 8).Value = 3

$q1.Cells.Item(9, 1).Value = 7
$q1.Cells.Item(9, 2).Value = "011356"
$q1.Cells.Item(9, 3).Value = "华泰柏瑞港股通时代机遇混合型证券投资基金C"
$q1.Cells.Item(9, 4).Value = "0.40"
$q1.Cells.Item(9, 5).Value = "90.93"
$q1.Cells.Item(9, 6).Value = "8.07"
$q1.Cells.Item(9, 7).Value = "0.0323"
$q1.Cells.Item(9, 8).Value = 3

$q1.Cells.Item(10, 1).Value = 8
$q1.Cells.Item(10, 2).Value = "011906"
$q1.Cells.Item(10, 3).Value = "安信价值启航混合型证券投资基金C"
$q1.Cells.Item(10, 4).Value = "0.42"
$q1.Cells.Item(10, 5).Value = "93.14"
$q1.Cells.Item(10, 6).Value = "3.21"
$q1.Cells.Item(10, 7).Value = "0.0135"
$q1.Cells.Item(10, 8).Value = 7

$q1.Cells.Item(11, 1).Value = 9
$q1.Cells.Item(11, 2).Value = "519602"
$q1.Cells.Item(11, 3).Value = "海富通大中华精选混合QDII"
$q1.Cells.Item(11, 4).Value = "0.11"
$q1.Cells.Item(11, 5).Value = "89.68"
$q1.Cells.Item(11, 6).Value = "5.45"
$q1.Cells.Item(11, 7).Value = "0.0060"
$q1.Cells.Item(11, 8).Value = 3

# --- Update the "总计" sheet: insert the new 2022-Q1 summary row at the top ---
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# The inserted blank row borrows row-1 (header) formatting; restore the plain
# "data row" look by pulling the index-column style from A3 and clearing the
# borrowed formatting on B2:D2 (matches the unstyled B/C/D data cells elsewhere).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.85

# Re-number the shifted index column (old A2:A6 holding 0..4 are now A3:A7).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

